# Generate Report for Handoff
# Updates the "b.md" row (row 3) on the Overview / zh-cn / de-de sheets to
# reflect that a new handoff has been generated for b.md.

$wb = $excel.ActiveWorkbook

function Get-HyperlinkForCell($sheet, $cellAddress) {
    $target = $sheet.Range($cellAddress).Address()
    for ($i = 1; $i -le $sheet.Hyperlinks.Count; $i++) {
        $hl = $sheet.Hyperlinks.Item($i)
        if ($hl.Range.Address() -eq $target) {
            return $hl
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Overview sheet: row 3 is the "b.md" row.
#   B3 (zh-cn status) / C3 (de-de status): "Handed back..." -> "Ready for handoff"
#   D3 (Latest Handoff Date): -> 2016-03-23 06:35:51
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-23 06:35:51"

# ---------------------------------------------------------------------------
# zh-cn sheet: row 3 is the "b.md" row.
#   C3 (Status): "Handed back..." -> "Ready for handoff"
#   D3 (Latest Handoff File): new handoff xlf name, with matching hyperlink
#   E3 (Latest Handoff Datetime): -> 2016-03-23 06:35:47
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"

$zhCnHandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c2fe69397ba09ccdde59353f85460ca5fbeea396/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhCnNewFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"

$zhCnLink = Get-HyperlinkForCell $wsZhCn "D3"
$zhCnLink.Address = $zhCnHandoffUrl
$zhCnLink.TextToDisplay = $zhCnNewFile

$wsZhCn.Range("E3").Value = "2016-03-23 06:35:47"

# ---------------------------------------------------------------------------
# de-de sheet: row 3 is the "b.md" row.
#   C3 (Status): "Handed back..." -> "Ready for handoff"
#   D3 (Latest Handoff File): new handoff xlf name, with matching hyperlink
#   E3 (Latest Handoff Datetime): -> 2016-03-23 06:35:51
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"

$deDeHandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/53fc40d024f618570b4418403271dec3f0faf965/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deDeNewFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"

$deDeLink = $wsDeDe.Hyperlinks.Item(6)
$deDeLink.Address = $deDeHandoffUrl
$deDeLink.TextToDisplay = $deDeNewFile

$wsDeDe.Range("E3").Value = "2016-03-23 06:35:51"

Write-Output "Handoff report updated for b.md"
